$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update PARENT_SITE_ID values in rows 2 and 3 from "1042001" to "1407001"
$ws.Range("A2").Value = "1407001"
$ws.Range("A3").Value = "1407001"
